$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "login_credentials" to "users"
$ws.Name = "users"

# New header cells for the added columns
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "phone"
$ws.Range("E1").Value = "age"

# Row 2 (user1 / pass1) - email address, to become a hyperlink
$ws.Range("C2").Value = "sid@test.com"
# Row 3 (user2 / pass2) - email address, to become a hyperlink
$ws.Range("C3").Value = "poudhan@test.com"

# Turn the two e-mail cells into mailto: hyperlinks (also applies the
# built-in "Hyperlink" style - underline + theme color)
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sid@test.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:poudhan@test.com")

# Phone & age, entered as quote-prefixed text (leading apostrophe keeps
# the numeric-looking string stored as text with the General format).
# Filled in column-by-column (both phone numbers, then both ages).
$ws.Range("D2").Value = "'2132138098"
$ws.Range("D3").Value = "'32840932"
$ws.Range("E2").Value = "'32"
$ws.Range("E3").Value = "'30"

# Size column D (phone numbers) the way Excel's AutoFit left it
$ws.Columns("D").ColumnWidth = 10.166666666666666

# Final selected cell as saved in the workbook
$ws.Range("E4").Select()
